$d = $word.ActiveDocument

$replacements = @(
    @{old="74÷8="; new="37÷6="},
    @{old="94÷8="; new="75÷3="},
    @{old="65÷2="; new="84÷2="},
    @{old="94÷6="; new="36÷2="},
    @{old="25÷3="; new="86÷4="},
    @{old="54÷7="; new="47÷5="},
    @{old="89÷6="; new="44÷3="},
    @{old="74÷5="; new="91÷8="},
    @{old="87÷3="; new="12÷4="},
    @{old="91÷6="; new="22÷8="},
    @{old="61÷6="; new="49÷6="},
    @{old="94÷5="; new="75÷4="},
    @{old="41÷6="; new="48÷2="},
    @{old="10÷7="; new="94÷7="},
    @{old="78÷3="; new="89÷2="},
    @{old="27÷5="; new="83÷8="},
    @{old="13÷5="; new="47÷5="},
    @{old="79÷8="; new="10÷8="},
    @{old="54÷3="; new="46÷5="},
    @{old="15÷8="; new="72÷9="},
    @{old="70÷5="; new="59÷3="},
    @{old="41÷8="; new="81÷9="},
    @{old="73÷6="; new="97÷5="},
    @{old="41÷7="; new="26÷5="},
    @{old="89÷4="; new="49÷3="}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
